$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A60").Value = "Matteo Azzolini"
$ws.Range("B60").Value = "MATTEO PILATI | Pinguini Trentini"
$ws.Range("C60").Value = "Mattia Baldessarini | Shark Attack"
$ws.Range("D60").Value = "ENRICO BORDIGNON | Pinguini Trentini"
$ws.Range("E60").Value = "Stefano  Galvagni | Clitoriders"
$ws.Range("F60").Value = "Davide  Bazzano  | iMontagna"
